# Update "想去人数" (want-to-go count) figures in column F across the
# relevant worksheets, reflecting the newly generated gh-pages output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2798
$ws1.Range("F5").Value  = 1555
$ws1.Range("F9").Value  = 126
$ws1.Range("F11").Value = 9456
$ws1.Range("F12").Value = 406
$ws1.Range("F13").Value = 2509
$ws1.Range("F15").Value = 267
$ws1.Range("F18").Value = 662
$ws1.Range("F20").Value = 1194
$ws1.Range("F22").Value = 2947
$ws1.Range("F23").Value = 2226
$ws1.Range("F25").Value = 1924
$ws1.Range("F26").Value = 1933
$ws1.Range("F28").Value = 1555
$ws1.Range("F29").Value = 300
$ws1.Range("F30").Value = 12
$ws1.Range("F34").Value = 339
$ws1.Range("F35").Value = 67
$ws1.Range("F38").Value = 19
$ws1.Range("F39").Value = 107
$ws1.Range("F40").Value = 1323
$ws1.Range("F41").Value = 111
$ws1.Range("F42").Value = 1459
$ws1.Range("F43").Value = 21
$ws1.Range("F44").Value = 331
$ws1.Range("F45").Value = 21
$ws1.Range("F46").Value = 203
$ws1.Range("F47").Value = 718

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 1
$ws2.Range("F10").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 2798
$ws4.Range("F4").Value  = 1555
$ws4.Range("F5").Value  = 1
$ws4.Range("F8").Value  = 126
$ws4.Range("F9").Value  = 9456
$ws4.Range("F10").Value = 406
$ws4.Range("F14").Value = 267
$ws4.Range("F16").Value = 662
$ws4.Range("F17").Value = 1194
$ws4.Range("F19").Value = 2947
$ws4.Range("F20").Value = 2226
$ws4.Range("F21").Value = 1924
$ws4.Range("F23").Value = 1555
$ws4.Range("F24").Value = 300
$ws4.Range("F25").Value = 12
$ws4.Range("F29").Value = 339
$ws4.Range("F30").Value = 67
$ws4.Range("F36").Value = 19
$ws4.Range("F37").Value = 107
$ws4.Range("F38").Value = 1323
$ws4.Range("F40").Value = 111
$ws4.Range("F41").Value = 1459
$ws4.Range("F42").Value = 21
$ws4.Range("F44").Value = 331
$ws4.Range("F45").Value = 21
$ws4.Range("F46").Value = 203
$ws4.Range("F47").Value = 718
$ws4.Range("F49").Value = 3
